# Update the "K" column (column G) values in Sheet1 of the workbook.
# These values were regenerated (K computed instead of the old "Strike#"
# derived values), so we write the new, recalculated values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 4
    12 = 3
    13 = 4
    14 = 2
    15 = 2
    16 = 7
    17 = 0
    18 = 5
    19 = 3
    20 = 1
    21 = 0
    22 = 5
    23 = 0
    24 = 2
    25 = 1
    26 = 2
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 6
    32 = 0
    33 = 1
    34 = 5
    35 = 3
    36 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
